# Updated cryptos list on Thu May 25 06:11:42 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking strings
# (e.g. "307.42") are stored as text, matching the original inline-string cells,
# instead of being auto-coerced into numbers by the Excel value parser.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '26.425.33'
$ws.Range('E2').Value = '  -1.78%  '
$ws.Range('D3').Value = '1.797.14'
$ws.Range('E3').Value = '  -2.10%  '
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '1.006'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').Value = '307.42'
$ws.Range('E6').Value = '  -1.08%  '
$ws.Range('D7').Value = '0.4567'
$ws.Range('E7').Value = '  -1.36%  '
$ws.Range('D8').Value = '0.3625'
$ws.Range('E8').Value = '  -1.12%  '
$ws.Range('D9').Value = '46.50'
$ws.Range('E9').Value = '  +1.26%  '
$ws.Range('D10').Value = '0.07115'
$ws.Range('E10').Value = '  -0.79%  '
$ws.Range('D11').Value = '0.8831'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').Value = '0.07822'
$ws.Range('E12').Value = '  -0.64%  '
$ws.Range('D13').Value = '19.53'
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '5.287'
$ws.Range('E14').Value = '  -1.12%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.719.25'
$ws.Range('E15').Value = '  -6.08%  '
$ws.Range('D16').Value = '6.332'
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('D17').Value = '85.14'
$ws.Range('E17').Value = '  -3.38%  '
$ws.Range('D18').Value = '1.007'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').Value = '0.000008588'
$ws.Range('E19').Value = '  -1.66%  '
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').Value = '26.462.92'
$ws.Range('E21').Value = '  -1.69%  '
$ws.Range('D22').Value = '14.28'
$ws.Range('E22').Value = '  -1.31%  '
$ws.Range('D23').Value = '4.994'
$ws.Range('E23').Value = '  -0.24%  '
$ws.Range('D24').Value = '10.51'
$ws.Range('E24').Value = '  +0.86%  '
$ws.Range('D25').Value = '1.945.25'
$ws.Range('E25').Value = '  -7.02%  '
$ws.Range('D26').Value = '1.979'
$ws.Range('E26').Value = '  -0.45%  '
$ws.Range('D27').Value = '152.58'
$ws.Range('E27').Value = '  +1.13%  '
$ws.Range('E28').Value = '  -1.89%  '
$ws.Range('D29').Value = '2.049'
$ws.Range('E29').Value = '  +3.30%  '
$ws.Range('D30').Value = '112.02'
$ws.Range('E30').Value = '  -1.42%  '
$ws.Range('D31').Value = '4.874'
$ws.Range('E31').Value = '  -1.33%  '
$ws.Range('D32').Value = '0.08657'
$ws.Range('E32').Value = '  -2.21%  '
$ws.Range('D33').Value = '3.054'
$ws.Range('E33').Value = '  -2.73%  '
$ws.Range('D34').Value = '4.456'
$ws.Range('E34').Value = '  -0.24%  '
$ws.Range('D35').Value = '0.7269'
$ws.Range('E35').Value = '  -4.31%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '2.715'
$ws.Range('E36').Value = '  +3.16%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '1.111'
$ws.Range('E37').Value = '  -1.54%  '
$ws.Range('E38').Value = '  +0.21%  '
$ws.Range('E39').Value = '  -1.01%  '
$ws.Range('D40').Value = '0.01946'
$ws.Range('E40').Value = '  +0.56%  '
$ws.Range('D41').Value = '0.05127'
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('D42').Value = '2.872'
$ws.Range('E42').Value = '  -2.01%  '
$ws.Range('D43').Value = '0.5192'
$ws.Range('E43').Value = '  +4.19%  '
$ws.Range('D44').Value = '6.902'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = '0.1529'
$ws.Range('E45').Value = '  -4.27%  '
$ws.Range('D46').Value = '8.034'
$ws.Range('E46').Value = '  -3.32%  '
$ws.Range('D47').Value = '0.4687'
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('B48').Value = 'PaxosStandard'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range('D48').Value = '1.007'
$ws.Range('E48').Value = '  -31.94%  '
$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').Value = '1.007'
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '9.885'
$ws.Range('E50').Value = '  -2.80%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = '100.47'
$ws.Range('E51').Value = '  -1.81%  '
# Restore the default (Normal) style on column D so no stray number-format
# style survives on the cells once the text values are locked in.
$priceRange.Style = "Normal"
